$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# ---------------------------------------------------------------------------
# 1. Flip the "Нет" -> "Да" cells that also get the new green highlight fill
#    (these were yellow-filled "Нет" cells; they become green-filled "Да").
#    Fill color = theme "Green, Accent 6, Lighter 60%" (RGB C5E0B4).
# ---------------------------------------------------------------------------
$targetCells = @("Q4","Z4","Z14","AB14","S15","Z15","AA15","AB15","AC15")
$greenColor = 11854021   # RGB(197, 224, 180) == &HC5E0B4 (BGR-packed for COM)

foreach ($addr in $targetCells) {
    $cell = $ws.Range($addr)
    $cell.Value = "Да"
    $cell.Interior.Color = $greenColor
}

# ---------------------------------------------------------------------------
# 2. Re-merge the header row ranges so the merged-cell list is rewritten with
#    the F1:F2 / A1:A2 / B1:B2 / C1:C2 / D1:D2 / E1:E2 block relocated ahead
#    of the rest (matches the saved file's new merge ordering).
# ---------------------------------------------------------------------------
$reorderRanges = @("M1:M2","N1:S1","T1:AF1","AG1:AH1","G1:G2","H1:H2","I1:I2","J1:J2","K1:K2","L1:L2")
foreach ($r in $reorderRanges) {
    $ws.Range($r).UnMerge()
    $ws.Range($r).Merge()
}

# ---------------------------------------------------------------------------
# 3. Update the view: scroll/select so the active cell becomes X15.
# ---------------------------------------------------------------------------
$ws.Range("X15").Select()
